$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 2 updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 189
$wsOff.Range("C2").Value = 140
$wsOff.Range("D2").Value = 31
$wsOff.Range("E2").Value = 8

# Sheet "DEF" - row 2 updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 149
$wsDef.Range("C2").Value = 108
$wsDef.Range("D2").Value = 24
$wsDef.Range("E2").Value = 13
$wsDef.Range("F2").Value = 1
